$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Last table row: bump the fixed row height (567 twips -> 680 twips,
#    and switch the height rule from the implicit "at least" to an exact
#    rule) to match the new grades-journal sizing.
# ---------------------------------------------------------------------------
$t = $d.Tables(1)
$lastRow = $t.Rows($t.Rows.Count)
$lastRow.HeightRule = 2   # wdRowHeightExactly
$lastRow.Height = 34      # 34pt == 680 twips

# ---------------------------------------------------------------------------
# 2) The "#kc" column-header cells: Word's proofer now flags "kc" as a
#    misspelled fragment, splitting the run into "#" + "kc" and wrapping
#    the latter in spell-check proofErr markers. Re-create that exact
#    paragraph shape (preserving each paragraph's own rsid/spacing) via
#    InsertXML so the replacement is a faithful, minimal content edit.
# ---------------------------------------------------------------------------
function Set-KcParagraph($cellIndex, $rsidRPr, $spacingBefore, $spacingLine, $hasBookmark) {
    $cell = $lastRow.Cells($cellIndex)
    $para = $cell.Range.Paragraphs(1)

    $tail = '<w:proofErr w:type="spellEnd"/>'
    if ($hasBookmark) {
        $tail = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + $tail
    }

    $body = '<w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:before="' + $spacingBefore + '" w:line="' + $spacingLine + '" w:lineRule="exact"/><w:rPr><w:sz w:val="16"/></w:rPr></w:pPr>' +
            '<w:r><w:rPr><w:sz w:val="16"/></w:rPr><w:t>#</w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:rPr><w:sz w:val="16"/></w:rPr><w:t>kc</w:t></w:r>' +
            $tail

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/temp.xml" pkg:contentType="application/xml">' +
           '<pkg:xmlData>' +
           '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
           'w:rsidR="008A2248" w:rsidRPr="' + $rsidRPr + '" w:rsidRDefault="00107261" w:rsidP="00DC45C8">' +
           $body +
           '</w:p></pkg:xmlData></pkg:part></pkg:package>'

    $para.Range.InsertXML($xml)
}

Set-KcParagraph 1  "00D8473E" "0" "140" $false
Set-KcParagraph 3  "00D8473E" "0" "138" $false
Set-KcParagraph 5  "00D8473E" "0" "139" $false
Set-KcParagraph 7  "002D0E3A" "1" "139" $false
Set-KcParagraph 9  "002D0E3A" "1" "135" $false
Set-KcParagraph 11 "002D0E3A" "2" "134" $true

# ---------------------------------------------------------------------------
# 3) "Normal Table" style: mark it as a quick style (adds <w:qFormat/>),
#    matching the refreshed template metadata.
# ---------------------------------------------------------------------------
$normalTableStyle = $d.Styles("Normal Table")
$normalTableStyle.QuickStyle = $true
